# edit.ps1 - Reproduce the commit "Changes to python document"
# Applies structural / content / formatting changes to the three worksheets
# of the Master_Classifications workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("body-based helpful")
$ws2 = $wb.Worksheets.Item("time-based helpful")
$ws3 = $wb.Worksheets.Item("otherwise")

# ---------------------------------------------------------------------------
# Sheet 1 "body-based helpful"
# ---------------------------------------------------------------------------
# Row 174 gets an explicit row height of 30 (auto height was 15 before).
$ws1.Rows(174).RowHeight = 30

# ---------------------------------------------------------------------------
# Sheet 2 "time-based helpful"
# ---------------------------------------------------------------------------
# G42 loses its red "needs review" fill and becomes the normal grey wrap-text
# style (same style already used elsewhere in column G, e.g. G72/G73).
$ws2.Range("G72").Copy()
$ws2.Range("G42").PasteSpecial(-4122)   # xlPasteFormats

# New cell B70 "okay are you ready" - same style already used in B2 (grey
# solid fill, wrap text, bordered).
$ws2.Range("B2").Copy()
$ws2.Range("B70").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("B70").Value = "okay are you ready"

# G78 & G79 lose their red "needs review" fill and become the plain wrap
# text style used by the rest of column G (e.g. G70/G71), and their
# content is replaced by what used to sit two rows further down (G80/G81),
# i.e. the two red flagged cells are effectively deleted and the column
# shifts up to fill the gap.
$g80 = $ws2.Range("G80").Value2
$g81 = $ws2.Range("G81").Value2
$g82 = $ws2.Range("G82").Value2
$g83 = $ws2.Range("G83").Value2
$g84 = $ws2.Range("G84").Value2
$g85 = $ws2.Range("G85").Value2
$g86 = $ws2.Range("G86").Value2

$ws2.Range("G70").Copy()
$ws2.Range("G78").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("G78").Value = $g80

$ws2.Range("G71").Copy()
$ws2.Range("G79").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("G79").Value = $g81

$ws2.Range("G80").Value = $g82
$ws2.Range("G81").Value = $g83
$ws2.Range("G82").Value = $g84
$ws2.Range("G83").Value = $g85

# G84 becomes the bold-bordered style (same as G77) carrying the text that
# used to be in G86.
$ws2.Range("G77").Copy()
$ws2.Range("G84").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("G84").Value = $g86

# G85 & G86 no longer hold any content - fully cleared (value + format).
$ws2.Range("G85").Clear()
$ws2.Range("G86").Clear()

# Row heights follow the content that now occupies each row (auto height
# is not recalculated headlessly, so set explicitly to match the new wrap
# text content).
$ws2.Rows(78).RowHeight = 75
$ws2.Rows(79).RowHeight = 105
$ws2.Rows(80).RowHeight = 135
$ws2.Rows(81).RowHeight = 45
$ws2.Rows(82).RowHeight = 60
$ws2.Rows(83).RowHeight = 45
$ws2.Rows(84).RowHeight = 120
$ws2.Rows(86).RowHeight = 15

# ---------------------------------------------------------------------------
# Sheet 3 "otherwise"
# ---------------------------------------------------------------------------
# Row 216 (the red-flagged "okay are you ready" row) is removed entirely;
# everything below shifts up by one row.
$ws3.Rows(216).Delete()

# The duplicate-value conditional formatting rule that used to sit on
# E219 now needs to point at E218 after the row shift above.
$fc = $ws3.Range("E219").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws3.Range("E218"))

# Two additional lines of dialogue are appended at the end of the sheet,
# plus one further empty row, using the same style as the existing last
# rows (column A unstyled, column E carrying style "3").
$ws3.Range("A300").Copy()
$ws3.Range("A301").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("A301").Value = "Calvin_3"

$ws3.Range("E300").Copy()
$ws3.Range("E301").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("E301").Value = "you ready for this Wade"

$ws3.Range("A301").Copy()
$ws3.Range("A302").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("A302").Value = "Calvin_3"

$ws3.Range("E301").Copy()
$ws3.Range("E302").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("E302").Value = "you ready Angie"

$ws3.Range("A302").Copy()
$ws3.Range("A303").PasteSpecial(-4122)  # xlPasteFormats
$ws3.Range("A303").Value = "Calvin_3"

# ---------------------------------------------------------------------------
# Active sheet / selections (cosmetic view state)
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("C4").Select()

$ws2.Activate()
$ws2.Range("F47").Select()

$ws3.Activate()
$ws3.Range("F212").Select()

$ws1.Activate()
